$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, bypassing Excel's automatic
# number/date inference, and without creating any new cell style.
# We do this by putting a formula that evaluates to a text string into
# the cell, then converting it to a literal value via copy/paste-values.
function Set-TextValue {
    param($addr, $text)
    $escaped = $text.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null  # xlPasteValues
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Row 1: new header cells I1:Y1 (same bold/bordered header style as the
# existing headers, copied from A1 so the existing style index is reused)
# ---------------------------------------------------------------------
$headers = @{
    "I1" = "Measurement Tool (Images)";
    "J1" = "Voltage (Images)";
    "K1" = "Expected Value (Images)";
    "L1" = "Polarisers (Images)";
    "M1" = "Cell Orientation";
    "N1" = "Polariser Number";
    "O1" = "State of Cell";
    "P1" = "Measurement Tool (Test B)";
    "Q1" = "Voltage (Test B)";
    "R1" = "Expected Value (Test B)";
    "S1" = "Test B Field 1";
    "T1" = "Test B Field 2";
    "U1" = "Measurement Tool (Test C)";
    "V1" = "Voltage (Test C)";
    "W1" = "Expected Value (Test C)";
    "X1" = "Test C Field 1";
    "Y1" = "Test C Field 2";
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

$ws.Range("A1").Copy() | Out-Null
$ws.Range("I1:Y1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Rows 2-4: materialize empty cells in the newly added columns I:Y
# (touching formatting without assigning a value creates an empty cell
# without introducing any new style record)
# ---------------------------------------------------------------------
$ws.Range("I2:Y4").Font.Bold = $false

# ---------------------------------------------------------------------
# Row 4: D4/E4 change from text "20"/"30" to real numbers 20/30
# ---------------------------------------------------------------------
$ws.Range("D4").Value = 20
$ws.Range("E4").Value = 30

# ---------------------------------------------------------------------
# Row 5: new data row
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "Transmittance"

# B5:E5 stay blank but present
$ws.Range("B5:E5").Font.Bold = $false

Set-TextValue "F5" "2025-01-29"
Set-TextValue "G5" "42421"

$ws.Range("H5").Value = 4

Set-TextValue "I5" "Tool A"
Set-TextValue "J5" "02"
Set-TextValue "K5" "12"
Set-TextValue "L5" "Yes"
Set-TextValue "M5" "45"
Set-TextValue "N5" "2"
Set-TextValue "O5" "ON"
Set-TextValue "P5" "Tool A"
Set-TextValue "Q5" "3"
Set-TextValue "R5" "4"
Set-TextValue "S5" "5"
Set-TextValue "T5" "6"
Set-TextValue "U5" "Tool A"
Set-TextValue "V5" "5"
Set-TextValue "W5" "5"
Set-TextValue "X5" "6"
Set-TextValue "Y5" "7"
